# week 5 - exam results updated
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New Hmw2 (column D) scores for rows that didn't have one yet ---
$ws.Range("D6").Value = 100
$ws.Range("D12").Value = 70
$ws.Range("D36").Value = 75
$ws.Range("D39").Value = 30
$ws.Range("D44").Value = 90
$ws.Range("D56").Value = 70
$ws.Range("D63").Value = 100
$ws.Range("D68").Value = 85
$ws.Range("D85").Value = 70

# --- New column E (second homework note) text values ---
# Write E86 ("*****") before E25 ("****") so the shared-string table picks
# up the same ordering as the authored workbook (index 355 = "*****",
# index 356 = "****").
$ws.Range("E86").Value = "*****"
$ws.Range("E25").Value = "****"

# --- View tweaks: zoom in on the sheet and move the selection ---
$excel.ActiveWindow.Zoom = 120
$ws.Range("C63:C64").Select()
